$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F
$ws.Range("F1").Value = "nummer6"

# New values for column F, rows 2-6 (matching the existing rows)
$ws.Range("F2").Value = 20220301
$ws.Range("F3").Value = 20220312
$ws.Range("F4").Value = 20220325
$ws.Range("F5").Value = 20220401
$ws.Range("F6").Value = 20220501

# New row 7 across columns A-F
$ws.Range("A7").Value = 20220601
$ws.Range("B7").Value = 20220602
$ws.Range("C7").Value = 20220601
$ws.Range("D7").Value = 20220604
$ws.Range("E7").Value = 20220605
$ws.Range("F7").Value = 20220605

# Match column width of the new column to the existing C:E block
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Update selection to match the post-edit state
$ws.Range("F3").Select() | Out-Null
